# Case_4_187 - 380 kV case: update computed power-flow values in pl_mw sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6793107954798927
$ws.Range("C2").Value = 0.1525804409479576
$ws.Range("D2").Value = 0.6034296570682898
$ws.Range("E2").Value = 0.2321810351232294
$ws.Range("G2").Value = 0.5951141136014826
$ws.Range("H2").Value = 0.7112117221102778
$ws.Range("I2").Value = 0.5071687533280347
$ws.Range("J2").Value = 0.1116288951982867
$ws.Range("M2").Value = 0.3808890920517243
$ws.Range("N2").Value = 1.175697335602941
$ws.Range("O2").Value = 2.593322707566358
$ws.Range("B3").Value = 0.6049434980060937
$ws.Range("C3").Value = 0.1336355259787751
$ws.Range("D3").Value = 0.5993216955707794
$ws.Range("E3").Value = 0.231797210348379
$ws.Range("G3").Value = 0.5929443060359603
$ws.Range("H3").Value = 0.7143397936445979
$ws.Range("I3").Value = 0.5122633347427055
$ws.Range("J3").Value = 0.1121671938242166
$ws.Range("M3").Value = 0.3584315496976913
$ws.Range("N3").Value = 1.184309024462664
$ws.Range("O3").Value = 2.594969703520292
$ws.Range("B4").Value = 0.5592629065957624
$ws.Range("C4").Value = 0.1219732711205381
$ws.Range("D4").Value = 0.5970895539580283
$ws.Range("E4").Value = 0.2316719298149721
$ws.Range("G4").Value = 0.5920428204128001
$ws.Range("H4").Value = 0.7166048894320625
$ws.Range("I4").Value = 0.515704910850733
$ws.Range("J4").Value = 0.1125481760014146
$ws.Range("M4").Value = 0.3447591368851235
$ws.Range("N4").Value = 1.190042829475999
$ws.Range("O4").Value = 2.597604444320723
$ws.Range("B5").Value = 0.5406441892657199
$ws.Range("C5").Value = 0.1172135310553415
$ws.Range("D5").Value = 0.5962530106407797
$ws.Range("E5").Value = 0.2316486600029144
$ws.Range("G5").Value = 0.591783657184024
$ws.Range("H5").Value = 0.7176145839856218
$ws.Range("I5").Value = 0.5171861381644014
$ws.Range("J5").Value = 0.1127161214875159
$ws.Range("M5").Value = 0.3392171827342523
$ws.Range("N5").Value = 1.192491716983213
$ws.Range("O5").Value = 2.599086161636507
$ws.Range("B6").Value = 0.5375523895051799
$ws.Range("C6").Value = 0.1164227482393301
$ws.Range("D6").Value = 0.5961185198638788
$ws.Range("E6").Value = 0.2316464748248883
$ws.Range("G6").Value = 0.5917471546556072
$ws.Range("H6").Value = 0.7177874774779553
$ws.Range("I6").Value = 0.5174368494996102
$ws.Range("J6").Value = 0.1127447753551571
$ws.Range("M6").Value = 0.3382987472696115
$ws.Range("N6").Value = 1.192905139962342
$ws.Range("O6").Value = 2.599356838466065
$ws.Range("B7").Value = 0.559011820645253
$ws.Range("C7").Value = 0.1219091087095308
$ws.Range("D7").Value = 0.5970779760294107
$ws.Range("E7").Value = 0.2316715034647423
$ws.Range("G7").Value = 0.5920388873180258
$ws.Range("H7").Value = 0.7166181556385567
$ws.Range("I7").Value = 0.5157245684099863
$ws.Range("J7").Value = 0.1125503895791873
$ws.Range("M7").Value = 0.3446842755858484
$ws.Range("N7").Value = 1.190075401112303
$ws.Range("O7").Value = 2.597622775410144
$ws.Range("B8").Value = 0.6536735522201695
$ws.Range("C8").Value = 0.1460546154737301
$ws.Range("D8").Value = 0.6019530709093885
$ws.Range("E8").Value = 0.2320257951262441
$ws.Range("G8").Value = 0.5942764769063729
$ws.Range("H8").Value = 0.7122188115469044
$ws.Range("I8").Value = 0.5088602592465783
$ws.Range("J8").Value = 0.1118040279527293
$ws.Range("M8").Value = 0.3731217273392815
$ws.Range("N8").Value = 1.178574107335066
$ws.Range("O8").Value = 2.59355351238267
$ws.Range("B9").Value = 0.8391116298465136
$ws.Range("C9").Value = 0.1931563308506838
$ws.Range("D9").Value = 0.6138118583105268
$ws.Range("E9").Value = 0.2335956804489143
$ws.Range("G9").Value = 0.6020893015473945
$ws.Range("H9").Value = 0.7063236476781043
$ws.Range("I9").Value = 0.497889978012239
$ws.Range("J9").Value = 0.110740783678235
$ws.Range("M9").Value = 0.4298006679067825
$ws.Range("N9").Value = 1.159555546921602
$ws.Range("O9").Value = 2.59846780313984
$ws.Range("B10").Value = 0.9751887899288931
$ws.Range("C10").Value = 0.2276015631105963
$ws.Range("D10").Value = 0.6239228467017597
$ws.Range("E10").Value = 0.2352820036168524
$ws.Range("G10").Value = 0.6099284898182447
$ws.Range("H10").Value = 0.7036571480849148
$ws.Range("I10").Value = 0.4913530884205422
$ws.Range("J10").Value = 0.1102037217639591
$ws.Range("M10").Value = 0.4719876161412202
$ws.Range("N10").Value = 1.147731636825547
$ws.Range("O10").Value = 2.609961466455957
$ws.Range("B11").Value = 1.037049143999241
$ws.Range("C11").Value = 0.2432349520425419
$ws.Range("D11").Value = 0.6288257039706195
$ws.Range("E11").Value = 0.2361647880559126
$ws.Range("G11").Value = 0.6139530441682552
$ws.Range("H11").Value = 0.7028054498587437
$ws.Range("I11").Value = 0.4887109996767762
$ws.Range("J11").Value = 0.1100124176124524
$ws.Range("M11").Value = 0.4912956523620053
$ws.Range("N11").Value = 1.142818047589046
$ws.Range("O11").Value = 2.616907162957972
$ws.Range("B12").Value = 1.06046697816555
$ws.Range("C12").Value = 0.2491495166373738
$ws.Range("D12").Value = 0.6307258243776914
$ws.Range("E12").Value = 0.2365156921231275
$ws.Range("G12").Value = 0.6155431321762705
$ws.Range("H12").Value = 0.7025348710240422
$ws.Range("I12").Value = 0.4877582735611057
$ws.Range("J12").Value = 0.109947598390697
$ws.Range("M12").Value = 0.4986236379316722
$ws.Range("N12").Value = 1.141024201342461
$ws.Range("O12").Value = 2.619784567834898
$ws.Range("B13").Value = 1.055423876063571
$ws.Range("C13").Value = 0.247875956879426
$ws.Range("D13").Value = 0.6303146657453453
$ws.Range("E13").Value = 0.236439380116984
$ws.Range("G13").Value = 0.6151977371301598
$ws.Range("H13").Value = 0.7025908351658217
$ws.Range("I13").Value = 0.4879613338281032
$ws.Range("J13").Value = 0.1099612192953678
$ws.Range("M13").Value = 0.4970446993473345
$ws.Range("N13").Value = 1.141407567264082
$ws.Range("O13").Value = 2.619153867371665
$ws.Range("B14").Value = 1.038975897980436
$ws.Range("C14").Value = 0.2437216579576216
$ws.Range("D14").Value = 0.62898115650367
$ws.Range("E14").Value = 0.2361933243154297
$ws.Range("G14").Value = 0.6140825365395131
$ws.Range("H14").Value = 0.7027821482873833
$ws.Range("I14").Value = 0.4886316602865932
$ws.Range("J14").Value = 0.1100069321034702
$ws.Range("M14").Value = 0.4918982025586089
$ws.Range("N14").Value = 1.142669127908263
$ws.Range("O14").Value = 2.617138932452775
$ws.Range("B15").Value = 1.028900041406246
$ws.Range("C15").Value = 0.2411763100406858
$ws.Range("D15").Value = 0.6281700074669345
$ws.Range("E15").Value = 0.2360447709073732
$ws.Range("G15").Value = 0.6134080535199473
$ws.Range("H15").Value = 0.7029060968947789
$ws.Range("I15").Value = 0.4890484792143432
$ws.Range("J15").Value = 0.1100359253554792
$ws.Range("M15").Value = 0.4887479561228858
$ws.Range("N15").Value = 1.143450570684792
$ws.Range("O15").Value = 2.615936930689486
$ws.Range("B16").Value = 0.971145115337265
$ws.Range("C16").Value = 0.2265791363261087
$ws.Range("D16").Value = 0.6236085282675958
$ws.Range("E16").Value = 0.2352266373871501
$ws.Range("G16").Value = 0.6096747140802847
$ws.Range("H16").Value = 0.7037200777166674
$ws.Range("I16").Value = 0.4915324354363548
$ws.Range("J16").Value = 0.110217290682602
$ws.Range("M16").Value = 0.4707281128591205
$ws.Range("N16").Value = 1.148062105954459
$ws.Range("O16").Value = 2.609542127842246
$ws.Range("B17").Value = 0.9357026824425247
$ws.Range("C17").Value = 0.2176148328560146
$ws.Range("D17").Value = 0.6208878230094115
$ws.Range("E17").Value = 0.2347543499135227
$ws.Range("G17").Value = 0.6075019626651539
$ws.Range("H17").Value = 0.7043119572695389
$ws.Range("I17").Value = 0.4931412544409923
$ws.Range("J17").Value = 0.1103421297677087
$ws.Range("M17").Value = 0.4597032208203373
$ws.Range("N17").Value = 1.151010224800473
$ws.Range("O17").Value = 2.606059141471633
$ws.Range("B18").Value = 0.9153132743749666
$ws.Range("C18").Value = 0.2124554392576101
$ws.Range("D18").Value = 0.6193515011343891
$ws.Range("E18").Value = 0.234493591031832
$ws.Range("G18").Value = 0.606295397856627
$ws.Range("H18").Value = 0.704686399399705
$ws.Range("I18").Value = 0.4940978083027048
$ws.Range("J18").Value = 0.1104189232081367
$ws.Range("M18").Value = 0.4533730275431722
$ws.Range("N18").Value = 1.152749692219068
$ws.Range("O18").Value = 2.604217430717227
$ws.Range("B19").Value = 0.9084091489100956
$ws.Range("C19").Value = 0.2107079885844882
$ws.Range("D19").Value = 0.6188362365407158
$ws.Range("E19").Value = 0.2344071732573028
$ws.Range("G19").Value = 0.6058942812559707
$ws.Range("H19").Value = 0.7048190205405547
$ws.Range("I19").Value = 0.4944270373636073
$ws.Range("J19").Value = 0.1104457810342332
$ws.Range("M19").Value = 0.4512316390946225
$ws.Range("N19").Value = 1.153346168647076
$ws.Range("O19").Value = 2.6036216069663
$ws.Range("B20").Value = 0.939475998826822
$ws.Range("C20").Value = 0.218569448497334
$ws.Range("D20").Value = 0.6211744920464639
$ws.Range("E20").Value = 0.2348034989604564
$ws.Range("G20").Value = 0.6077287893984504
$ws.Range("H20").Value = 0.7042454310160764
$ws.Range("I20").Value = 0.4929667627505552
$ws.Range("J20").Value = 0.1103283240529471
$ws.Range("M20").Value = 0.4608757002071968
$ws.Range("N20").Value = 1.15069186117902
$ws.Range("O20").Value = 2.606413183062699
$ws.Range("B21").Value = 1.043807277581209
$ws.Range("C21").Value = 0.244942026800004
$ws.Range("D21").Value = 0.6293716603287862
$ws.Range("E21").Value = 0.2362651461635004
$ws.Range("G21").Value = 0.6144083035000421
$ws.Range("H21").Value = 0.7027245454425355
$ws.Range("I21").Value = 0.4884334717822618
$ws.Range("J21").Value = 0.1099932982353593
$ws.Range("M21").Value = 0.4934094098872208
$ws.Range("N21").Value = 1.142296763913365
$ws.Range("O21").Value = 2.61772405586504
$ws.Range("B22").Value = 1.111950406928543
$ws.Range("C22").Value = 0.2621460665789925
$ws.Range("D22").Value = 0.6349825735778438
$ws.Range("E22").Value = 0.2373172396510057
$ws.Range("G22").Value = 0.6191589657653509
$ws.Range("H22").Value = 0.702033302916476
$ws.Range("I22").Value = 0.4857492207092662
$ws.Range("J22").Value = 0.1098187756499662
$ws.Range("M22").Value = 0.5147677887678128
$ws.Range("N22").Value = 1.137199552047072
$ws.Range("O22").Value = 2.626557489149491
$ws.Range("B23").Value = 1.075585543782438
$ws.Range("C23").Value = 0.2529669689736238
$ws.Range("D23").Value = 0.6319647538213076
$ws.Range("E23").Value = 0.2367468649275821
$ws.Range("G23").Value = 0.6165881517429739
$ws.Range("H23").Value = 0.7023745357102342
$ws.Range("I23").Value = 0.487156340244816
$ws.Range("J23").Value = 0.1099078553754005
$ws.Range("M23").Value = 0.5033597840945703
$ws.Range("N23").Value = 1.139884414542784
$ws.Range("O23").Value = 2.621710960265801
$ws.Range("B24").Value = 0.9377701231348965
$ws.Range("C24").Value = 0.2181378844868789
$ws.Range("D24").Value = 0.6210448022248443
$ws.Range("E24").Value = 0.2347812451418001
$ws.Range("G24").Value = 0.6076261084126457
$ws.Range("H24").Value = 0.704275401118295
$ws.Range("I24").Value = 0.4930455519484127
$ws.Range("J24").Value = 0.1103345499745849
$ws.Range("M24").Value = 0.4603455968188399
$ws.Range("N24").Value = 1.150835654580533
$ws.Range("O24").Value = 2.606252620293475
$ws.Range("B25").Value = 0.7889712117043359
$ws.Range("C25").Value = 0.180441553139616
$ws.Range("D25").Value = 0.6103579597892264
$ws.Range("E25").Value = 0.233077356090579
$ws.Range("G25").Value = 0.5996079050923555
$ws.Range("H25").Value = 0.7076260553265854
$ws.Range("I25").Value = 0.5005906137641141
$ws.Range("J25").Value = 0.1109855545498029
$ws.Range("M25").Value = 0.4143709680216787
$ws.Range("N25").Value = 1.164322743471445
$ws.Range("O25").Value = 2.595755619893168

Write-Host "Updated 264 cells for 380 kV case"
